# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (cloned from "2021-Q4" to keep the same
#   layout / cell formatting) right after "2021-Q4" and before "总计".
# - Populate it with the new holdings data.
# - Insert a corresponding new row at the top of the "总计" (summary) sheet's
#   data, shifting the older rows down by one and renumbering the index
#   column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the "2022-Q1" worksheet, cloned from "2021-Q4" so it keeps the
#    exact same column layout / styles, then place it right before "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Row 1 headers already match the template ("基金代码", "基金名称", ...,
# "持有市值(亿元)", "仓位排名") so only the data row needs updating.
# The code/name/size/position figures are stored as plain text (not
# numbers) in the source data, so a leading apostrophe forces text entry;
# resetting the style back to Normal afterwards drops the "quote prefix"
# flag Excel would otherwise remember, leaving a plain unstyled text cell.
$newSheet.Range("B2").Value = "'378006"
$newSheet.Range("C2").Value = "'上投摩根全球新兴市场混合(QDII)"
$newSheet.Range("D2").Value = "'0.46"
$newSheet.Range("E2").Value = "'88.99"
$newSheet.Range("F2").Value = "'1.97"
$newSheet.Range("G2").Value = "'0.0091"
$newSheet.Range("B2:G2").Style = "Normal"
$newSheet.Range("H2").Value = 8

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q1" row into the "总计" (summary) sheet, pushing
#    the existing rows down one and renumbering the A-column index.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Shift the existing data rows (2-6) down to (3-7) preserving their
# formatting/types.
$summary.Range("A2:D6").Copy($summary.Range("A3:D7"))

# Renumber the shifted rows' index column.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# Write the new summary row for 2022-Q1.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q1"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01
